$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.183.37"
$ws.Range("E2").Value = "  -0.12%  "
$ws.Range("D3").Value = "1.904.47"
$ws.Range("E3").Value = "  +0.45%  "
$ws.Range("E4").Value = "  -0.45%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "253.96"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.25%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.700"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.42%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "41.77"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.55%  "
$ws.Range("E9").Value = "  +2.07%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "52.79"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.55%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0761"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.70%  "
$ws.Range("E12").Value = "  -0.37%  "
$ws.Range("E13").Value = "  +3.33%  "
$ws.Range("D14").Value = "2.181.72"
$ws.Range("E14").Value = "  +0.52%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.738"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.54%  "
$ws.Range("E16").Value = "  +4.82%  "
$ws.Range("D17").Value = "1.886.80"
$ws.Range("E17").Value = "  -0.50%  "
$ws.Range("D18").Value = "35.167.82"
$ws.Range("E18").Value = "  -0.18%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "73.88"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.68%  "
$ws.Range("D20").Value = "0.0₃0845"
$ws.Range("E20").Value = "  +3.30%  "
$ws.Range("E21").Value = "  +1.11%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "13.08"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.52%  "
$ws.Range("E23").Value = "  +5.81%  "
$ws.Range("E24").Value = "  -0.39%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.45"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.30%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.34"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.30%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "168.05"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.21%  "
$ws.Range("E28").Value = "  -0.32%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.57"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.44%  "
$ws.Range("E30").Value = "  -0.01%  "
$ws.Range("D31").Value = "4.128.22"
$ws.Range("E31").Value = "  -0.35%  "
$ws.Range("E32").Value = "  +6.48%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.35"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.77%  "
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0601"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.02%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.63"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +9.38%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.25"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.58%  "
$ws.Range("E37").Value = "  -0.45%  "
$ws.Range("E38").Value = "  -6.51%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.02"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.07%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "98.72"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +9.63%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "17.15"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.56%  "
$ws.Range("E42").Value = "  +4.47%  "
$ws.Range("E43").Value = "  +1.92%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0658"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.89%  "
$ws.Range("E45").Value = "  -0.02%  "
$ws.Range("D46").Value = "1.307.87"
$ws.Range("E46").Value = "  -2.96%  "
$ws.Range("E47").Value = "  +0.09%  "
$ws.Range("E48").Value = "  -1.28%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "12.29"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.92%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.61"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.40%  "
$ws.Range("E51").Value = "  +7.45%  "
